$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows containing "Robin" and "Reehan" (rows 3 and 4), shifting "Rahul" up to row 3
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Move the active selection to D11, matching the final saved selection state
$ws.Range("D11").Select()
